$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted right before the current
# row 337, which pushes the existing rows 337:413 down to 338:414 and
# extends the used range to A1:R414.
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row 337 with the new record
# (Acelga, Primera, Provincia de Cautín, week of 2022-11-24).
$ws.Cells.Item(337, 1).Value = 10
$ws.Cells.Item(337, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(337, 3).Value = "La Araucanía"
$ws.Cells.Item(337, 4).Value = 44889
$ws.Cells.Item(337, 5).Value = 9
$ws.Cells.Item(337, 6).Value = 100112009
$ws.Cells.Item(337, 7).Value = "Acelga"
$ws.Cells.Item(337, 8).Value = "Sin especificar"
$ws.Cells.Item(337, 9).Value = "Primera"
$ws.Cells.Item(337, 10).Value = 85
$ws.Cells.Item(337, 11).Value = 9000
$ws.Cells.Item(337, 12).Value = 9000
$ws.Cells.Item(337, 13).Value = 9000
$ws.Cells.Item(337, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(337, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(337, 16).Value = 750
$ws.Cells.Item(337, 17).Value = 12
$ws.Cells.Item(337, 18).Value = "Hortaliza"
